$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$lo = $ws2.ListObjects.Item(1)

# --- Expand the XML-mapped table on sheet "Foglio2" (bug 62055 resize fix) ---

# Insert a new column so the table grows from 3 to 4 columns; this pushes the
# existing "SchemaRef" / "Namespace" columns one column to the right and gives
# the newly inserted cells the formatting of the column to their left, which
# mirrors how Excel grows a table when XML data with extra columns is imported.
$ws2.Range("D1").EntireColumn.Insert()

# The full-row column insert also shifted the two single (unmapped) XML cells
# that live outside the table (H2/H3) one column over (to I2/I3). Put them
# back where they belong.
$ws2.Cells.Item(2, 9).Copy($ws2.Cells.Item(2, 8))
$ws2.Cells.Item(3, 9).Copy($ws2.Cells.Item(3, 8))
$ws2.Cells.Item(2, 9).Clear()
$ws2.Cells.Item(3, 9).Clear()

# Grow the table/list object definition itself to include the new column.
$lo.Resize($ws2.Range("C5:F9"))

# Give the new column its header, "Unmapped Column" (a column with no XML
# mapping), and re-touch the other header cells so the table metadata picks
# up their (shifted) names correctly.
$ws2.Cells.Item(5, 4).Value2 = "Unmapped Column"
$ws2.Cells.Item(5, 5).Value2 = $ws2.Cells.Item(5, 5).Value2
$ws2.Cells.Item(5, 6).Value2 = $ws2.Cells.Item(5, 6).Value2

# Add the note cell just below/left of the table that documents the expected
# behaviour; it should be cleared automatically when the table expands.
$note = $ws2.Cells.Item(11, 4)
$note.Value2 = "[This text should be cleared if the table expands]"
$note.Font.Italic = $true

# --- Update sheet selections / active sheet to match the saved workbook state ---

$ws1.Activate()
$ws1.Range("D1").Select()

$ws3.Activate()
$ws3.Range("C36").Select()

$ws2.Activate()
$ws2.Range("A1").Select()
